# Update column G ("K" = strikeouts) values on Sheet1 to reflect the
# regenerated save_data (was using a "Strike#"-style metric, now uses
# actual strikeout counts K).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$kValues = @{
    2  = 11
    3  = 5
    4  = 7
    5  = 7
    6  = 8
    7  = 7
    8  = 9
    9  = 3
    10 = 3
    11 = 4
    12 = 6
    13 = 6
    14 = 4
    15 = 1
    16 = 4
    17 = 4
    18 = 7
    19 = 0
    20 = 4
    21 = 10
    22 = 4
    23 = 6
    24 = 3
    25 = 6
    26 = 7
    27 = 10
    28 = 6
    29 = 3
    30 = 3
    31 = 4
    33 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
